$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a "last changed" date that was bulk
# updated from 2023-09-06 (Excel serial 45175) to 2023-09-08 (serial
# 45177) for every data row in the sheet.
$oldSerial = 45175
$newDate = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0).Date.AddDays(45177)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq $oldSerial) {
        $cell.Value = $newDate
    }
}
